$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 259.66666
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H76").Value = 166672850
$ws.Range("I76").Value = 333339700
$ws.Range("K76").Value = 333339700
$ws.Range("M76").Value = -333339385
$ws.Range("H79").Value = 166672850
$ws.Range("I79").Value = 333339700
$ws.Range("K79").Value = 333339700
$ws.Range("M79").Value = -333338608
$ws.Range("H82").Value = 1496
$ws.Range("I82").Value = 1496
$ws.Range("K82").Value = 4488
$ws.Range("M82").Value = -4082
$ws.Range("H85").Value = 1496
$ws.Range("I85").Value = 1496
$ws.Range("K85").Value = 4488
$ws.Range("M85").Value = -3084
$ws.Range("H92").Value = 1355
$ws.Range("I92").Value = 1355
$ws.Range("K92").Value = 1355
$ws.Range("M92").Value = -107
$ws.Range("H138").Value = 3435.8333
$ws.Range("I138").Value = 797.95654
$ws.Range("K138").Value = 2393.86962
$ws.Range("M138").Value = 2746.13038
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8945.857
$ws.Range("I74").Value = 11026.2
$ws.Range("J74").Value = 3745
$ws.Range("K74").Value = 11026.2
$ws.Range("L74").Value = 3745
$ws.Range("M74").Value = -10152.2
$ws.Range("N74").Value = -5493
$ws.Range("H77").Value = 8945.857
$ws.Range("I77").Value = 11026.2
$ws.Range("J77").Value = 3745
$ws.Range("K77").Value = 55131
$ws.Range("L77").Value = 18725
$ws.Range("M77").Value = -50763
$ws.Range("N77").Value = -27461
$ws.Range("H97").Value = 1362.1666
$ws.Range("I97").Value = 1218.5
$ws.Range("K97").Value = 1218.5
$ws.Range("M97").Value = -722.5
$ws.Range("H122").Value = 6950.3
$ws.Range("J122").Value = 6984
$ws.Range("L122").Value = 20952
$ws.Range("N122").Value = -25852
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 50003720
$ws.Range("I16").Value = 71431096
$ws.Range("J16").Value = 6499
$ws.Range("K16").Value = 71431096
$ws.Range("L16").Value = 6499
$ws.Range("M16").Value = -71430809
$ws.Range("N16").Value = -7073
$ws.Range("H22").Value = 478128.9
$ws.Range("I22").Value = 993192.9399999999
$ws.Range("K22").Value = 993192.9399999999
$ws.Range("M22").Value = -992842.9399999999
$ws.Range("H31").Value = 5441.294
$ws.Range("I31").Value = 1836.875
$ws.Range("J31").Value = 8645.223
$ws.Range("K31").Value = 1836.875
$ws.Range("L31").Value = 8645.223
$ws.Range("M31").Value = -1541.875
$ws.Range("N31").Value = -9235.223
$ws.Range("H34").Value = 5441.294
$ws.Range("I34").Value = 1836.875
$ws.Range("J34").Value = 8645.223
$ws.Range("K34").Value = 1836.875
$ws.Range("L34").Value = 8645.223
$ws.Range("M34").Value = -1634.875
$ws.Range("N34").Value = -9049.223
$ws.Range("H58").Value = 71444984
$ws.Range("I58").Value = 111123020
$ws.Range("J58").Value = 24499
$ws.Range("K58").Value = 111123020
$ws.Range("L58").Value = 24499
$ws.Range("M58").Value = -111122817
$ws.Range("N58").Value = -24905
$ws.Range("H62").Value = 8061.3335
$ws.Range("I62").Value = 10175.833
$ws.Range("J62").Value = 3832.3333
$ws.Range("K62").Value = 10175.833
$ws.Range("L62").Value = 3832.3333
$ws.Range("M62").Value = -9551.833000000001
$ws.Range("N62").Value = -5080.3333
$ws.Range("H65").Value = 8061.3335
$ws.Range("I65").Value = 10175.833
$ws.Range("J65").Value = 3832.3333
$ws.Range("K65").Value = 50879.165
$ws.Range("L65").Value = 19161.6665
$ws.Range("M65").Value = -47759.165
$ws.Range("N65").Value = -25401.6665
$ws.Range("H107").Value = 1158
$ws.Range("I107").Value = 899
$ws.Range("J107").Value = 1201.1666
$ws.Range("K107").Value = 899
$ws.Range("L107").Value = 1201.1666
$ws.Range("M107").Value = 1021
$ws.Range("N107").Value = -5041.1666
$ws.Range("H113").Value = 50003720
$ws.Range("I113").Value = 71431096
$ws.Range("J113").Value = 6499
$ws.Range("K113").Value = 71431096
$ws.Range("L113").Value = 6499
$ws.Range("M113").Value = -71428926
$ws.Range("N113").Value = -10839
$ws.Range("H122").Value = 5738.769
$ws.Range("I122").Value = 5335.6
$ws.Range("J122").Value = 5990.75
$ws.Range("K122").Value = 16006.8
$ws.Range("L122").Value = 17972.25
$ws.Range("M122").Value = -13556.8
$ws.Range("N122").Value = -22872.25
$ws.Range("H134").Value = 125020620
$ws.Range("I134").Value = 166684830
$ws.Range("J134").Value = 27998.5
$ws.Range("K134").Value = 500054490
$ws.Range("L134").Value = 83995.5
$ws.Range("M134").Value = -500051955
$ws.Range("N134").Value = -89065.5
$ws.Range("H136").Value = 71444984
$ws.Range("I136").Value = 111123020
$ws.Range("J136").Value = 24499
$ws.Range("K136").Value = 333369060
$ws.Range("L136").Value = 73497
$ws.Range("M136").Value = -333366510
$ws.Range("N136").Value = -78597
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3930.7693
$ws.Range("I122").Value = 2164.3333
$ws.Range("K122").Value = 6492.999899999999
$ws.Range("M122").Value = -4042.999899999999
$ws.Range("H133").Value = 299999.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2284.7083
$ws.Range("I68").Value = 2304.2
$ws.Range("J68").Value = 2187.25
$ws.Range("K68").Value = 2304.2
$ws.Range("L68").Value = 2187.25
$ws.Range("M68").Value = -1555.2
$ws.Range("N68").Value = -3685.25
$ws.Range("H71").Value = 2284.7083
$ws.Range("I71").Value = 2304.2
$ws.Range("J71").Value = 2187.25
$ws.Range("K71").Value = 11521
$ws.Range("L71").Value = 10936.25
$ws.Range("M71").Value = -7777
$ws.Range("N71").Value = -18424.25
$ws.Range("H122").Value = 3434
$ws.Range("I122").Value = 3434
$ws.Range("K122").Value = 10302
$ws.Range("M122").Value = -7852
$ws.Range("H125").Value = 49613
$ws.Range("J125").Value = 49613
$ws.Range("L125").Value = 49613
$ws.Range("N125").Value = -59453
$ws.Range("H132").Value = 2885.1428
$ws.Range("I132").Value = 1348
$ws.Range("K132").Value = 4044
$ws.Range("M132").Value = -1514
$ws.Range("H136").Value = 51728940
$ws.Range("I136").Value = 23813396
$ws.Range("J136").Value = 125007250
$ws.Range("K136").Value = 71440188
$ws.Range("L136").Value = 375021750
$ws.Range("M136").Value = -71437638
$ws.Range("N136").Value = -375026850
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5159.4
$ws.Range("I122").Value = 5699.25
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 17097.75
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -14647.75
$ws.Range("N122").Value = -13900
$ws.Range("H136").Value = 26322308
$ws.Range("I136").Value = 29412888
$ws.Range("K136").Value = 88238664
$ws.Range("M136").Value = -88236114
